# Cập nhật tên bài tập: "Bài tập 1" -> "Bài tập 12.1" on the slide titled
# "Bài tập 1" (slide 31), splitting the title into three runs the way
# PowerPoint does when a user edits text in place.

$p = $ppt.ActivePresentation

$target = $null
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    if ($slide.Shapes.HasTitle) {
        $titleShape = $slide.Shapes.Title
        if ($titleShape.TextFrame.TextRange.Text -eq "Bài tập 1") {
            $target = $slide
            break
        }
    }
}

$title = $target.Shapes.Title
$tr = $title.TextFrame.TextRange

# Replace the trailing "1" with "12.1", then split "tập " off from "Bài ".
$tr.Characters(9, 1).Text = "12.1"
$tr.Characters(5, 4).Text = "tập "
